$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: price changed from "160 р." to "157 р."
$ws.Range("B2").Value = "157 р."

# New row: B3 holds the text "30" (stored as text, not a number)
$ws.Range("B3").Value = "'30"
$ws.Range("B3").ClearFormats()

# Update the active selection to match B2
$ws.Range("B2").Select() | Out-Null
